$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B32").Value = ".Integração com sistema de polícia para carros roubados"
$fnt = $ws.Range("B32").Font
$fnt.Color = 255
$fnt.Bold = $true
